$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet1 ("工作表1"): insert two blocks of new rows (13 rows before the old
# row 6, and 2 more rows before what becomes row 37) and populate them.
# ---------------------------------------------------------------------------
$ws1.Range("A6:A18").EntireRow.Insert()
$ws1.Range("A35:A36").EntireRow.Insert()

$sheet1NewRows = @(
  @(6, "百齡橋進城", @(7299,5266,6666,8115)),
  @(7, "百齡橋出城", @(5284,6298,6637,6666)),
  @(8, "大直橋進城", @(2623,2584,8607,6813)),
  @(9, "大直橋出城", @(3503,2802,6813,8607)),
  @(10, "承德入城(劍潭路)", @(2053,1367,8289,15667)),
  @(11, "承德離城(劍潭路)", @(1383,1735,15667,8289)),
  @(12, "中山橋入城(劍潭路)", @(4836,3648,8291,10265)),
  @(13, "中山橋入城(通河路)", @(3012,1643,6675,8343)),
  @(14, "中山北離城(通河路)", @(1008,2029,8342,6675)),
  @(15, "新生高圓山端上匝道", @(3371,2459,6675,27271)),
  @(16, "新生高圓山端下匝道", @(2039,1799,27271,6675)),
  @(17, "復興北地下道南向", @(2104,2242,6814,6760)),
  @(18, "復興北地下道北向", @(2124,1974,6760,6814)),
  @(35, "環河北(-敦煌)南向", @(2610,1973,6652,8138)),
  @(36, "環河北(敦煌-)北向", @(2508,2685,8137,8136))
)

foreach ($item in $sheet1NewRows) {
  $r = $item[0]
  $name = $item[1]
  $vals = $item[2]
  $ws1.Cells.Item($r, 1).Value = $name
  for ($c = 0; $c -lt $vals.Length; $c++) {
    $ws1.Cells.Item($r, $c + 2).Value = $vals[$c]
  }
}

# Restore the selection shown in the saved file (single cell E37).
$ws1.Range("E37").Select()

# ---------------------------------------------------------------------------
# Sheet2 ("工作表2"): append a new row of data (row 14) and set its page
# setup (A4 portrait).
# ---------------------------------------------------------------------------
$ws2.Cells.Item(14, 1).Value = "中山北離城(劍潭路)"
$ws2.Cells.Item(14, 2).Value = 3431
$ws2.Cells.Item(14, 3).Value = 4494
$ws2.Cells.Item(14, 4).Value = 10266
$ws2.Cells.Item(14, 5).Value = 8291
$ws2.Cells.Item(14, 6).Value = 8334
$ws2.Cells.Item(14, 7).Value = 8341

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Restore sheet2's selection (A6:A13) with A6 active.
$ws2.Range("A6:A13").Select()

# Re-activate sheet1 so it is the sheet shown when the workbook is opened,
# and nudge the application window position to match the saved view.
$ws1.Activate()
$ws1.Range("E37").Select()
$excel.ActiveWindow.Left = 6288
$excel.ActiveWindow.Top = 2400
